$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new "Nucleo Virtual Comm Port" label for B4 and B5 (previously blank/red
# "unused" cells) - the UART pins were moved to the Nucleo's own virtual comm port.
$ws.Range("B4").Value = "Nucleo Virtual Comm Port"
$ws.Range("B5").Value = "Nucleo Virtual Comm Port"

# Those cells are no longer "do not use" (red) - recolour them the same grey used for
# other already-available/free pins.
$ws.Range("B4:B5").Interior.Color = 13421772

# UART_RX / UART_TX pins freed up - clear their purpose cells and drop the highlight fill
$ws.Range("F7").Value = ""
$ws.Range("F12").Value = ""
$ws.Range("F7").Interior.Pattern = -4142
$ws.Range("F12").Interior.Pattern = -4142

# Re-select F7 to match the final view state
$ws.Range("F7").Select()
